$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.577.51'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '3.088.38'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("D5").Value = '591.90'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").Value = '155.47'
$ws.Range("E6").Value = '  +7.22%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +3.04%  '
$ws.Range("D9").Value = '3.083.62'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").Value = '5.85'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = '37.60'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").Value = '3.600.10'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '7.21'
$ws.Range("E17").Value = '  -1.40%  '
$ws.Range("D18").Value = '63.544.28'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Value = '3.083.67'
$ws.Range("E19").Value = '  -1.72%  '
$ws.Range("D20").Value = '475.87'
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D21").Value = '14.70'
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").Value = '0.720'
$ws.Range("E22").Value = '  -1.60%  '
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("E24").Value = '  +4.17%  '
$ws.Range("D25").Value = '12.95'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '81.31'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '10.01'
$ws.Range("E27").Value = '  +2.46%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("E34").Value = '  -1.66%  '
$ws.Range("D35").Value = '0.0₃0852'
$ws.Range("E35").Value = '  +1.15%  '
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Value = '3.40'
$ws.Range("E37").Value = '  +6.66%  '
$ws.Range("D38").Value = '6.14'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  -2.99%  '
$ws.Range("D40").Value = '9.36'
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").Value = '50.82'
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D42").Value = '444.02'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("E44").Value = '  -1.90%  '
$ws.Range("D45").Value = '40.13'
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("E46").Value = '  +3.81%  '
$ws.Range("D47").Value = '2.806.86'
$ws.Range("E47").Value = '  -3.52%  '
$ws.Range("D48").Value = '131.88'
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").Value = '25.52'
$ws.Range("E49").Value = '  +5.59%  '
$ws.Range("E51").Value = '  +1.26%  '
